$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.945.32'
$ws.Range('E2').Value = '  +2.66%  '
$ws.Range('D3').Value = '2.960.62'
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  +0.05%  '
$c = $ws.Range('D5')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '595.12'
$c.Style = $origStyle
$ws.Range('E5').Value = '  +0.22%  '
$c = $ws.Range('D6')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '146.82'
$c.Style = $origStyle
$ws.Range('E6').Value = '  +1.23%  '
$c = $ws.Range('D7')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = $origStyle
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '2.960.05'
$ws.Range('E8').Value = '  +0.98%  '
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('E10').Value = '  +3.58%  '
$ws.Range('E11').Value = '  +5.77%  '
$ws.Range('E12').Value = '  +0.67%  '
$ws.Range('E13').Value = '  +5.99%  '
$c = $ws.Range('D14')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '33.13'
$c.Style = $origStyle
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').Value = '3.451.82'
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('D17').Value = '62.850.89'
$ws.Range('E17').Value = '  +2.70%  '
$c = $ws.Range('D18')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '6.74'
$c.Style = $origStyle
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = '2.957.17'
$ws.Range('E19').Value = '  +0.83%  '
$c = $ws.Range('D20')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '443.05'
$c.Style = $origStyle
$ws.Range('E20').Value = '  +2.13%  '
$c = $ws.Range('D21')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '13.48'
$c.Style = $origStyle
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('E22').Value = '  -1.78%  '
$c = $ws.Range('D23')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.12'
$c.Style = $origStyle
$ws.Range('E23').Value = '  -0.43%  '
$c = $ws.Range('D24')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '81.57'
$c.Style = $origStyle
$ws.Range('E24').Value = '  -0.58%  '
$c = $ws.Range('D25')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '11.21'
$c.Style = $origStyle
$ws.Range('E25').Value = '  +1.85%  '
$c = $ws.Range('D26')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '11.92'
$c.Style = $origStyle
$ws.Range('E26').Value = '  +0.30%  '
$c = $ws.Range('D27')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.14'
$c.Style = $origStyle
$ws.Range('E27').Value = '  -3.52%  '
$ws.Range('E28').Value = '  -0.01%  '
$c = $ws.Range('D29')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.32'
$c.Style = $origStyle
$ws.Range('E29').Value = '  +4.36%  '
$ws.Range('E30').Value = '  +0.67%  '
$c = $ws.Range('D31')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.17'
$c.Style = $origStyle
$ws.Range('E31').Value = '  -3.06%  '
$ws.Range('D32').Value = '0.0₃0967'
$ws.Range('E32').Value = '  +8.56%  '
$c = $ws.Range('D33')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '26.52'
$c.Style = $origStyle
$ws.Range('E33').Value = '  -1.02%  '
$ws.Range('E34').Value = '  -2.04%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  -1.44%  '
$c = $ws.Range('D37')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.12'
$c.Style = $origStyle
$ws.Range('E37').Value = '  +3.47%  '
$c = $ws.Range('D38')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.66'
$c.Style = $origStyle
$ws.Range('E38').Value = '  -0.41%  '
$ws.Range('E39').Value = '  +1.48%  '
$c = $ws.Range('D40')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '49.55'
$c.Style = $origStyle
$ws.Range('E40').Value = '  -0.66%  '
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('E42').Value = '  -5.19%  '
$c = $ws.Range('D43')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.281'
$c.Style = $origStyle
$ws.Range('E43').Value = '  -0.78%  '
$c = $ws.Range('D44')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '40.21'
$c.Style = $origStyle
$ws.Range('E44').Value = '  -5.33%  '
$ws.Range('D45').Value = '2.717.81'
$ws.Range('E45').Value = '  +0.49%  '
$c = $ws.Range('D46')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '134.37'
$c.Style = $origStyle
$ws.Range('E46').Value = '  +0.59%  '
$ws.Range('E47').Value = '  -3.04%  '
$c = $ws.Range('D48')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '362.89'
$c.Style = $origStyle
$ws.Range('E48').Value = '  -2.30%  '
$ws.Range('E50').Value = '  -0.72%  '
$c = $ws.Range('D51')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '22.86'
$c.Style = $origStyle
$ws.Range('E51').Value = '  -4.59%  '
